$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text format so values like "1.00" keep their
#     exact digits instead of being auto-coerced into numbers by Excel. ---
$priceCells = [ordered]@{
    'D2' = '57.740.16'
    'D3' = '2.441.08'
    'D4' = '1.00'
    'D5' = '524.62'
    'D6' = '130.07'
    'D7' = '0.998'
    'D9' = '2.442.36'
    'D12' = '4.93'
    'D14' = '2.873.64'
    'D15' = '57.673.62'
    'D16' = '21.61'
    'D18' = '2.435.21'
    'D19' = '10.29'
    'D20' = '4.12'
    'D21' = '312.48'
    'D22' = '6.07'
    'D23' = '0.998'
    'D24' = '64.78'
    'D25' = '0.402'
    'D27' = '0.157'
    'D28' = '7.22'
    'D29' = '173.62'
    'D31' = '1.69'
    'D36' = '17.76'
    'D38' = '3.75'
    'D39' = '36.36'
    'D41' = '0.788'
    'D43' = '263.79'
    'D46' = '0.0923'
    'D47' = '121.32'
    'D50' = '16.94'
    'D51' = '16.31'
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
    $cell.Style = "Normal"
}

# --- Volume(1h) column (E): plain text assignment (values already contain
#     spaces/% so Excel stores them as text without any extra coercion). ---
$volumeCells = [ordered]@{
    'E2' = '  -1.19%  '
    'E3' = '  -2.83%  '
    'E4' = '  +0.13%  '
    'E5' = '  +0.60%  '
    'E6' = '  -1.71%  '
    'E7' = '  -0.10%  '
    'E8' = '  +0.76%  '
    'E9' = '  -2.81%  '
    'E10' = '  -0.08%  '
    'E11' = '  -2.57%  '
    'E12' = '  -4.08%  '
    'E13' = '  -2.93%  '
    'E14' = '  -2.78%  '
    'E15' = '  -1.21%  '
    'E16' = '  -1.94%  '
    'E17' = '  -1.68%  '
    'E18' = '  -2.99%  '
    'E19' = '  -3.04%  '
    'E20' = '  -0.88%  '
    'E21' = '  -2.75%  '
    'E22' = '  -1.22%  '
    'E23' = '  -0.11%  '
    'E24' = '  +0.72%  '
    'E25' = '  -0.06%  '
    'E26' = '  +0.40%  '
    'E27' = '  -2.20%  '
    'E28' = '  -2.05%  '
    'E29' = '  +3.71%  '
    'E30' = '  -2.69%  '
    'E31' = '  -1.19%  '
    'E32' = '  -2.64%  '
    'E33' = '  -4.81%  '
    'E34' = '  -0.01%  '
    'E35' = '  -0.21%  '
    'E36' = '  -1.64%  '
    'E37' = '  -5.11%  '
    'E38' = '  -5.07%  '
    'E39' = '  +0.64%  '
    'E40' = '  -1.21%  '
    'E41' = '  +1.73%  '
    'E42' = '  -2.31%  '
    'E43' = '  -5.16%  '
    'E44' = '  -2.48%  '
    'E45' = '  -5.07%  '
    'E46' = '  +0.51%  '
    'E47' = '  -1.79%  '
    'E48' = '  -1.73%  '
    'E49' = '  -1.47%  '
    'E50' = '  -4.20%  '
    'E51' = '  -2.88%  '
}
foreach ($addr in $volumeCells.Keys) {
    $ws.Range($addr).Value = $volumeCells[$addr]
}
